# "Generate Report for Handback" — refresh the handback status report
# with the latest xliff-generation / handback timestamps and priority
# codes for the "ae138424-fc94-4960-a027-583fc84be83f.md" row (row 3)
# and the "f23d1d0e-c448-4e74-8061-ec51cc996248.md" row (row 5), which
# share the same values across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-09-04 14:18:46"
$wsOverview.Range("G5").Value = "2016-09-04 14:18:46"

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-09-04 14:18:42"
$wsZhCn.Range("K3").Value = "2016-09-04 14:18:59"
$wsZhCn.Range("E5").Value = "mt"
$wsZhCn.Range("H5").Value = "2016-09-04 14:18:42"
$wsZhCn.Range("K5").Value = "2016-09-04 14:18:59"

# --- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-09-04 14:18:46"
$wsDeDe.Range("K3").Value = "2016-09-04 14:19:10"
$wsDeDe.Range("E5").Value = "mt"
$wsDeDe.Range("H5").Value = "2016-09-04 14:18:46"
$wsDeDe.Range("K5").Value = "2016-09-04 14:19:10"
